$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with the reordered column labels
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "living_rooms_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "kitchens_2"
$ws.Range("F1").Value = "bedrooms_2"

# Update the 0/1 indicator matrix (rows 2-7) to match the new column order
$matrix = @(
    @(0, 1, 0, 0, 0, 0),
    @(0, 0, 0, 0, 1, 0),
    @(0, 0, 0, 1, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(1, 0, 0, 0, 0, 0),
    @(0, 0, 1, 0, 0, 0)
)

for ($i = 0; $i -lt $matrix.Length; $i++) {
    $rowValues = $matrix[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowValues[$j]
    }
}
